$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "Abdo"
$ws.Range("B9").Value = "IR HOLDER"
$ws.Range("C9").Value = "Rahman"

$ws.Columns.Item(4).EntireColumn.AutoFit()
